# Update cryptocurrency price/volume data per upstream GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.528.56"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").Value = "1.604.37"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("E6").Value = "  +7.15%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.89"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("E10").Value = "  +2.67%  "
$ws.Range("E11").Value = "  +2.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0909"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").Value = "1.833.70"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("D14").Value = "1.597.39"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.537"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.76%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "29.541.06"
$ws.Range("E16").Value = "  +2.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.75%  "
$ws.Range("E20").Value = "  +3.42%  "
$ws.Range("D21").Value = "0.0₃0692"
$ws.Range("E21").Value = "  +1.92%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("E27").Value = "  +5.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +2.75%  "
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("D35").Value = "1.408.66"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("E37").Value = "  +4.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.77%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("E41").Value = "  +3.82%  "
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0491"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "53.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +24.77%  "
$ws.Range("E45").Value = "  +3.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").Value = "1.744.81"
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.856"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "86.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.84%  "
